$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLO Stats-this session")

# Row 2: Raymond -> Fish
$ws.Range("A2").Value = "Fish"
$ws.Range("B2").Value = 39.46
$ws.Range("C2").Value = 19.46
$ws.Range("D2").Value = -20
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.857
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.3
$ws.Range("J2").Value = 0.429
$ws.Range("K2").Value = 0.357
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 28.4
$ws.Range("P2").Value = 3.8
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 0.833

# Row 3: Fish -> Cedric
$ws.Range("A3").Value = "Cedric"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 25.37
$ws.Range("D3").Value = 5.37
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.643
$ws.Range("G3").Value = 0.143
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.483
$ws.Range("J3").Value = 0.429
$ws.Range("K3").Value = 0.071
$ws.Range("L3").Value = 1.4
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 3.6
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 14
$ws.Range("R3").Value = 0.167
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T3").Value = "05/27/21"
$ws.Range("T3").Style = $ws.Range("T2").Style
